# Actualizacion automatica del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in results for rows that previously had no outcome recorded yet
$ws.Cells.Item(9, 7).Value = "Fallo"
$ws.Cells.Item(9, 8).Value = -1

$ws.Cells.Item(15, 7).Value = "Fallo"
$ws.Cells.Item(15, 8).Value = -1

# Append the newest tracked matches (rows 17-22), results still pending
$newRows = @(
    @(14716618, "2025-09-20", "Lorenzo Musetti", "Dino Prižmić", "Gana Lorenzo Musetti", 1.44),
    @(14687077, "2025-09-20", "Dalibor Svrcina", "Zhizhen Zhang", "Gana Dalibor Svrcina", 2.38),
    @(14655135, "2025-09-19", "Rafael Jodar", "Daniel Masur", "Gana Daniel Masur", 3.5),
    @(14718196, "2025-09-19", "Max Wiskandt", "Matthew Summers", "Gana Matthew Summers", 2.1),
    @(14718704, "2025-09-19", "Branko Djuric", "Michel Hopp", "Gana Michel Hopp", 3),
    @(14717098, "2025-09-19", "Gilles Arnaud Bailly", "Lorenzo Sciahbasi", "Gana Lorenzo Sciahbasi", 3.4)
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Keep "fecha" as plain text (e.g. "2025-09-20") instead of being
    # auto-converted to a date serial number.
    $ws.Cells.Item($r, 2).NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # resultado / profit are not known yet for these freshly added matches -
    # materialize the cells (matching the existing blank G/H cells used
    # elsewhere in the sheet for pending results) without leaving them
    # completely absent from the row.
    $ws.Cells.Item($r, 7).Style = $ws.Cells.Item(16, 7).Style
    $ws.Cells.Item($r, 8).Style = $ws.Cells.Item(16, 7).Style
}
